$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("db_corsi")

# New row 5 values
$courseName = "CORSO IN MATERIA DI PREVENZIONE DELLA CORRUZIONE E DELL'ILLEGALITA"
$contenuti  = "Definizioni operative del fenomeno corruttivo , Quadro normativo: art. 318 c.p. e art. 2635 c.c., Tipologie di corruzione: propria, impropria, ambientale e tra privati, Pene e circostanze aggravanti, Responsabilita individuali e aziendali, Policy interna anti-corruzione: obiettivi e ambito di applicazione, Procedure di segnalazione interna ed esterna, con garanzia di anonimato, Best practice etiche e strategie di compliance per il contrasto della corruzione"
$savePath   = "Anti Corruzione"

$ws.Range("A5").Value = $courseName
$ws.Range("B5").Value = $courseName
$ws.Range("C5").Value = $contenuti
$ws.Range("D5").Value = 8
$ws.Range("E5").Value = $savePath

# Match formatting of the rest of the sheet (row above) for A/B/D/E
$ws.Range("A5:E5").Font.Name = "Calibri"
$ws.Range("A5:E5").Font.Size = 11
$ws.Range("A5:E5").WrapText = $false

# Description cell gets its own wrap-text style
$ws.Range("C5").Font.Name = "Calibri"
$ws.Range("C5").Font.Size = 11
$ws.Range("C5").WrapText = $true

$ws.Rows.Item(5).RowHeight = 32.5

$ws.Range("C5").Select()
